# Rename the original sheet to "income" and add a new "income_by_category"
# sheet that summarizes total orders and income per category.

$wb = $excel.ActiveWorkbook

# 1. Rename the existing sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "income"

# 2. Add a new worksheet for the category summary, placed right after "income".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "income_by_category"

# Header labels (category / number_of_orders_by_item / income).
$ws2.Range("A1").Value = "category"
$ws2.Range("B1").Value = "number_of_orders_by_item"
$ws2.Range("C1").Value = "income"

# Data rows - totals by category.
$ws2.Range("A2").Value = "Asian"
$ws2.Range("B2").Value = 3470
$ws2.Range("A3").Value = "Italian"
$ws2.Range("B3").Value = 2948
$ws2.Range("A4").Value = "Mexican"
$ws2.Range("B4").Value = 2945
$ws2.Range("A5").Value = "American"
$ws2.Range("B5").Value = 2734

# Income column holds pre-formatted currency text (e.g. "$62,286.50"), so
# force text interpretation while assigning it, then drop the temporary
# number format again.
$ws2.Range("C2:C5").NumberFormat = "@"
$ws2.Range("C2").Value = "$62,286.50"
$ws2.Range("C3").Value = "$42,746.00"
$ws2.Range("C4").Value = "$38,137.75"
$ws2.Range("C5").Value = "$19,138.00"
$ws2.Range("C2:C5").ClearFormats()

# Give the header row the same bold/centered/bordered look used on "income".
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep the original sheet active/selected, matching the unedited source state.
$ws1.Activate()
